# Direct commit: put the text "This is change" into cell A1 of Sheet1,
# then select the whole column (as happens when a user clicks the column
# header to resize it) and set the column width to best-fit the content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "This is change"

$col = $ws.Columns.Item(1)
$col.ColumnWidth = 12.5
$col.Select() | Out-Null
